# Add season-record columns (Wins / Losses / Ties) to the player table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers (AD1:AF1) ---------------------------------------------------
# Copy the format of the last existing header cell (AC1, style s="1":
# bold, centered, bordered) onto the three new header cells so they match
# the look of the rest of the header row, then set their text.
$ws.Range("AC1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# --- Data (AD2:AF51) ------------------------------------------------------
# Every player row gets the same team season record: 82 wins, 80 losses,
# 0 ties.
for ($r = 2; $r -le 51; $r++) {
    $ws.Cells.Item($r, 30).Value = 82  # AD
    $ws.Cells.Item($r, 31).Value = 80  # AE
    $ws.Cells.Item($r, 32).Value = 0   # AF
}
